# energy carrier tax rate and mark-up updates for historic and future
# scenario simulations. presentation of the RokiG JF of 03.06 is added.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Sheet1")
$wsNote = $wb.Worksheets.Item("note")

# --- Update the tax-rate / mark-up values on Sheet1 ---
# Row 2 (id_energy_carrier = 1) and row 5 (id_energy_carrier = 6): all the
# yearly values (columns E:AT, years 2010-2051) are replaced with a flat
# 0.19 rate.
$wsData.Range("E2:AT2").Value = 0.19
$wsData.Range("E5:AT5").Value = 0.19

# --- Selection / active sheet bookkeeping ---
# Sheet1 becomes the active sheet/tab, with cell E7 selected.
$wsData.Activate() | Out-Null
$wsData.Range("E7").Select() | Out-Null

# The "note" sheet's selection stays at A3 (unchanged), but it is no longer
# the active tab.
$wsNote.Range("A3").Select() | Out-Null

# Selecting on "note" re-activates it as a side effect, so re-activate
# Sheet1 (with E7 selected) to make it the workbook's active tab on save.
$wsData.Activate() | Out-Null
$wsData.Range("E7").Select() | Out-Null
